# Revise audio controller board to use demultiplexers instead of relay
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Insert a new row before the old Relay row (row 6), pushing rows 6-9 (and the
# totals row) down by one. The new row will hold the R1,R2 resistor entry and
# the old Relay row (now row 7) is repurposed for the new U1,U2 demultiplexer.
$ws.Rows.Item(6).Insert()

# --- Row 2: C1 Ceramic Capacitor (value/part/price changed) ---
$ws.Cells.Item(2, 3).Value = "10µF"
$ws.Cells.Item(2, 4).Value = "CL21A106KOQNNNE"
$ws.Cells.Item(2, 9).Value = 0.1

# --- Row 3: was D1 Rectifier Diode -> C2 Ceramic Capacitor ---
$ws.Cells.Item(3, 1).Value = "C2"
$ws.Cells.Item(3, 2).Value = "Ceramic Capacitor"
$ws.Cells.Item(3, 3).Value = "1µF"
$ws.Cells.Item(3, 4).Value = "CL21B105KAFNNNE"
$ws.Cells.Item(3, 7).Value = "'0805"
$ws.Cells.Item(3, 9).Value = 0.1

# --- Row 4: was Q1 Basic FET P-Channel -> C3 Ceramic Capacitor ---
$ws.Cells.Item(4, 1).Value = "C3"
$ws.Cells.Item(4, 2).Value = "Ceramic Capacitor"
$ws.Cells.Item(4, 3).Value = "47nF"
$ws.Cells.Item(4, 7).Value = "'0805"
$ws.Cells.Item(4, 9).Value = 0.1

# --- Row 5: was R1 Resistor -> C4, C5 Ceramic Capacitor ---
$ws.Cells.Item(5, 1).Value = "C4, C5"
$ws.Cells.Item(5, 2).Value = "Ceramic Capacitor"
$ws.Cells.Item(5, 3).Value = "100nF"
$ws.Cells.Item(5, 4).Value = "CL21B104KACNNNC"
$ws.Cells.Item(5, 7).Value = "'0805"
$ws.Cells.Item(5, 8).Value = 2
$ws.Cells.Item(5, 9).Value = 0.1

# --- Row 6 (newly inserted): R1, R2 Resistor 10kOhm ---
$ws.Cells.Item(6, 2).Value = "Resistor"
$ws.Cells.Item(6, 3).Value = "10kΩ; ±5%; 1/4w"
$ws.Cells.Item(6, 4).Value = "ERJ-PA3F1002V"
$ws.Cells.Item(6, 1).Value = "R1, R2"
$ws.Cells.Item(6, 6).Value = "SMD"
$ws.Cells.Item(6, 7).Value = "'0603"
$ws.Cells.Item(6, 8).Value = 2
$ws.Cells.Item(6, 9).Value = 0.32

# --- Row 7: was RL1 Relay -> U1, U2 1:2 demultiplexer ---
$ws.Cells.Item(7, 1).Value = "U1, U2"
$ws.Cells.Item(7, 2).Value = "1:2 demultiplexer"
$ws.Cells.Item(7, 3).ClearContents()
$ws.Cells.Item(7, 4).Value = "SN74LVC1G18DBVR"
$ws.Cells.Item(7, 7).Value = "SOT-23-6"
$ws.Cells.Item(7, 8).Value = 2
$ws.Cells.Item(7, 9).Value = 0.82

# Rows 8, 9, 10 (J1, J2, J3/J4) keep their same content, just shifted down -
# no value changes needed there.

# Fix totals formulas to include the new row range (2:10 instead of 2:9)
$ws.Cells.Item(11, 8).Formula = "=SUM(H2:H10)"
$ws.Cells.Item(11, 9).Formula = "=SUM(I2:I10)"

# Update the table range to include the new row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:J11"))

# Part number for C3 was the last detail filled in (fixed up after the rest
# of the table was populated).
$ws.Cells.Item(4, 4).Value = "CL21B473KBCNNNC"

$ws.Range("J16").Select()
